$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in row 17 with new tracked entry
$ws.Range("A17").Value = 7.12
$ws.Range("B17").Value = 0.375
$ws.Range("C17").Value = 0.58333333333333337
$ws.Range("E17").Value = "5hr"
$ws.Range("F17").Value = "Get advices about function A filtering better solution include method and try to start Html layout"

# Match time formatting used by the rest of the Time IN / Time OUT columns
$ws.Range("B17:C17").NumberFormat = $ws.Range("B10:C10").NumberFormat

# Move the active selection to F18 (as recorded in the sheet view)
$ws.Range("F18").Select()
